# Update metrics table (columns B:Q, rows 2:26) with new values.
# All rows share the same new metric values per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.9999845286516351,
    0.9991389381953254,
    0.9999348357207491,
    0.9999896231028684,
    0.9999566782581165,
    0.00001444182458493536,
    0.000803763398418822,
    0.00007492135611277314,
    0.000006949727592248796,
    0.00004093554185251097,
    0.0001971186733916371,
    0.003800240069381849,
    1.00001125188972,
    0.003962024156716867,
    136.2907641523688,
    205.7666861698563
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
